$d = $word.ActiveDocument

# Confirm we can find the closing sentence of the "Day 3" paragraph -
# this is the last paragraph in the document, right where the new
# "Day 4" notes need to be appended.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "have more time to write the reports.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # The trailing "_GoBack" bookmark currently sits right after that
    # sentence (before the paragraph's end mark). Remove it here - it
    # gets re-created at the end of the newly appended content below,
    # matching where Word leaves it after the most recent edit.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    # New paragraphs, exactly as authored: a blank spacer paragraph, a
    # bold "Day 4" heading paragraph, and a paragraph with three runs
    # of body text (kept as separate runs, the way the original
    # content was produced), all using the same hr-HR language
    # formatting used throughout the rest of the notes.
    $newParagraphsXml = @'
<w:p><w:pPr><w:rPr><w:lang w:val="hr-HR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:lang w:val="hr-HR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="hr-HR"/></w:rPr><w:t>Day 4</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="hr-HR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="hr-HR"/></w:rPr><w:t xml:space="preserve">We are trying to work on our refactoring and today we will implement transactions. Some of the group members will start writing the report and we have decided on doing more Unit testing for the game to do the Quality Assurance. </w:t></w:r><w:r><w:rPr><w:lang w:val="hr-HR"/></w:rPr><w:t xml:space="preserve">We have managed to implement for the Web client to show the correct </w:t></w:r><w:r><w:rPr><w:lang w:val="hr-HR"/></w:rPr><w:t>and incorrect answers.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

    $insertionPoint = $d.Content
    $insertionPoint.Collapse(0)
    $insertionPoint.InsertXML($newParagraphsXml)
}
